# "Create reference data" rationalisation
#
# The "Updated" sheet used to carry two description columns
# (DESCRIPTION_OLD / DESCRIPTION_NEW) that both held the same kind of
# text. They are rationalised down to a single DESCRIPTION column that
# uses the new wording ("Anti-dumping / anti-subsidy") and the redundant
# third column is removed.

$wb = $excel.ActiveWorkbook
$wsUpdated = $wb.Worksheets.Item("Updated")
$wsNew     = $wb.Worksheets.Item("New")

# --- Rationalise the "Updated" sheet ---------------------------------

# The sheet used to carry both the old wording (column B,
# DESCRIPTION_OLD) and the new wording (column C, DESCRIPTION_NEW).
# Keep only the new wording, rename its header simply to DESCRIPTION,
# and drop the now-redundant old column.
$wsUpdated.Range("C1").Value = "DESCRIPTION"
$wsUpdated.Columns.Item(2).Delete()

# Match the refreshed column widths.
$wsUpdated.Columns.Item(1).ColumnWidth = 29.830729166666668
$wsUpdated.Columns.Item(2).ColumnWidth = 56.666666666666664

# --- Cosmetic selection / active sheet updates ------------------------

$wsNew.Activate()
$wsNew.Range("A2:B4").Select()
$wsNew.Columns.Item(1).ColumnWidth = 25.666666666666668
$wsNew.Columns.Item(2).ColumnWidth = 39.666666666666664

$wsUpdated.Activate()
$wsUpdated.Range("B2").Select()
